$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Add the new "description" header in column M of row 1
$ws.Cells.Item(1, 13).Value = "description"

# Select M1, matching the diff's updated selection
$ws.Range("M1").Select()
